$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row values
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of connector words ("de"/"del"/"el" -> "De"/"Del"/"El")
$ws.Range("A18").Value = "Ciudad De México"
$ws.Range("A22").Value = "Estado De México"
$ws.Range("B22").Value = "Ecatepec De Morelos"
$ws.Range("B27").Value = "San Luis De La Paz"
$ws.Range("B29").Value = "Acapulco De Juárez"
$ws.Range("B30").Value = "Atlamajalcingo Del Monte"
$ws.Range("B33").Value = "Mártir De Cuilapan"
$ws.Range("B48").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B49").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B55").Value = "Huehuetlán El Chico"
$ws.Range("B59").Value = "Amealco De Bonfil"
$ws.Range("B61").Value = "Jalpan De Serra"
$ws.Range("B62").Value = "Landa De Matamoros"
$ws.Range("B64").Value = "San Juan Del Río"
$ws.Range("B76").Value = "Soledad De Doblado"

# Remove trailing metadata/footer rows (86-90), shrinking the used range to A1:D84
$ws.Range("A86:D90").EntireRow.Delete()
